$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11: add response_collected for Ruan's r1 re-submission
$ws.Range("E11").Value = "master_worker_response_tracke_ruan_r1_Sep-30-2023.csv"

# Insert a new row at position 12 (pushes current row 12 and below down by one)
$ws.Rows.Item(12).Insert()

# New row 12: Ruan's second round submission
$ws.Range("A12").Value = "pairwise"
$ws.Range("B12").Value = "Sept-30-2023"
$ws.Range("C12").Value = "RD"
$ws.Range("D12").Value = "all_submitted_tracker_ruan_r2_Sep-30-2023.csv"

# New row 15 (appended at the end): SB HIT launch
$ws.Range("A15").Value = "paiewise"
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "Oct-1-2023"
$ws.Range("C15").Value = "SB"
$ws.Range("D15").Value = "all_submitted_tracker_SB_Oct-01-2023.csv"

$ws.Range("D15").Select()
